$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H40").Value = 4140.2974
$ws.Range("I40").Value = 3950
$ws.Range("J40").Value = 4231.64
$ws.Range("K40").Value = 3950
$ws.Range("L40").Value = 4231.64
$ws.Range("M40").Value = -3775
$ws.Range("N40").Value = -4581.64

$ws.Range("H70").Value = 6985.6665
$ws.Range("I70").Value = 1307.5
$ws.Range("J70").Value = 8608
$ws.Range("K70").Value = 3922.5
$ws.Range("L70").Value = 25824
$ws.Range("M70").Value = -3652.5
$ws.Range("N70").Value = -26364

$ws.Range("H73").Value = 6985.6665
$ws.Range("I73").Value = 1307.5
$ws.Range("J73").Value = 8608
$ws.Range("K73").Value = 3922.5
$ws.Range("L73").Value = 25824
$ws.Range("M73").Value = -2986.5
$ws.Range("N73").Value = -27696

$ws.Range("H116").Value = 8717.941999999999
$ws.Range("I116").Value = 6514.375
$ws.Range("K116").Value = 6514.375
$ws.Range("M116").Value = -3072.375

$ws.Range("H125").Value = 1616.4546
$ws.Range("J125").Value = 1383.2222
$ws.Range("L125").Value = 12448.9998
$ws.Range("N125").Value = -17368.9998

$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 5486.3403
$ws.Range("I61").Value = 4582.4287
$ws.Range("K61").Value = 4582.4287
$ws.Range("M61").Value = -4370.4287

$ws.Range("H96").Value = 27365.334
$ws.Range("J96").Value = 27365.334
$ws.Range("L96").Value = 27365.334
$ws.Range("N96").Value = -32857.334

$ws.Range("H132").Value = 1954.7833
$ws.Range("I132").Value = 1556.431
$ws.Range("J132").Value = 13507
$ws.Range("K132").Value = 4669.293
$ws.Range("L132").Value = 40521
$ws.Range("M132").Value = -2139.293
$ws.Range("N132").Value = -45581

$ws.Range("H136").Value = 5486.3403
$ws.Range("I136").Value = 4582.4287
$ws.Range("K136").Value = 13747.2861
$ws.Range("M136").Value = -11197.2861

$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 3119.25
$ws.Range("I86").Value = 2448.5
$ws.Range("J86").Value = 5131.5
$ws.Range("K86").Value = 2448.5
$ws.Range("L86").Value = 5131.5
$ws.Range("M86").Value = -1325.5
$ws.Range("N86").Value = -7377.5

$ws.Range("H89").Value = 3119.25
$ws.Range("I89").Value = 2448.5
$ws.Range("J89").Value = 5131.5
$ws.Range("K89").Value = 12242.5
$ws.Range("L89").Value = 25657.5
$ws.Range("M89").Value = -6626.5
$ws.Range("N89").Value = -36889.5

$ws.Range("H92").Value = 13923
$ws.Range("J92").Value = 13923
$ws.Range("L92").Value = 13923
$ws.Range("N92").Value = -18915

$ws.Range("H105").Value = 22135.857
$ws.Range("I105").Value = 24435
$ws.Range("K105").Value = 24435
$ws.Range("M105").Value = -22688

$ws.Range("H134").Value = 3859.1428
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 15014
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 45042
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -50112

$ws = $wb.Worksheets.Item(4)
$ws.Range("H42").Value = 4999.5
$ws.Range("I42").Value = 4999.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 4999.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -4406.5
$ws.Range("N42").ClearContents()

$ws.Range("H94").Value = 2782.25
$ws.Range("I94").Value = 1794.8
$ws.Range("J94").Value = 3487.5715
$ws.Range("K94").Value = 1794.8
$ws.Range("L94").Value = 3487.5715
$ws.Range("M94").Value = -1343.8
$ws.Range("N94").Value = -4389.5715

$ws.Range("H134").Value = 2379.889
$ws.Range("I134").Value = 1464.0625
$ws.Range("J134").Value = 9706.5
$ws.Range("K134").Value = 4392.1875
$ws.Range("L134").Value = 29119.5
$ws.Range("M134").Value = -1857.1875
$ws.Range("N134").Value = -34189.5

$ws = $wb.Worksheets.Item(5)
$ws.Range("H93").Value = 22513.25
$ws.Range("J93").Value = 22513.25
$ws.Range("L93").Value = 67539.75
$ws.Range("N93").Value = -71283.75

$ws.Range("H94").Value = 9693.666999999999
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H98").Value = 520.6
$ws.Range("J98").Value = 583.3333
$ws.Range("L98").Value = 1749.9999
$ws.Range("N98").Value = -4745.9999

$ws.Range("H112").Value = 83342560
$ws.Range("J112").Value = 15665
$ws.Range("L112").Value = 46995
$ws.Range("N112").Value = -49211

$ws.Range("H114").Value = 82.59999999999999
$ws.Range("I114").Value = 250
$ws.Range("J114").Value = 40.75
$ws.Range("K114").Value = 750
$ws.Range("L114").Value = 122.25
$ws.Range("M114").Value = 2504
$ws.Range("N114").Value = -6630.25

$ws = $wb.Worksheets.Item(6)
$ws.Range("H34").Value = 20000
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20536

$ws.Range("H64").Value = 24499.5
$ws.Range("J64").Value = 24499.5
$ws.Range("L64").Value = 24499.5
$ws.Range("N64").Value = -24995.5

$ws.Range("H67").Value = 24499.5
$ws.Range("J67").Value = 24499.5
$ws.Range("L67").Value = 24499.5
$ws.Range("N67").Value = -26215.5

$ws.Range("H70").Value = 7322.45
$ws.Range("I70").Value = 5379.5293
$ws.Range("K70").Value = 5379.5293
$ws.Range("M70").Value = -5109.5293

$ws.Range("H73").Value = 7322.45
$ws.Range("I73").Value = 5379.5293
$ws.Range("K73").Value = 5379.5293
$ws.Range("M73").Value = -4443.5293

$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630

$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184

$ws.Range("H126").Value = 6727.8184
$ws.Range("I126").Value = 4247.25
$ws.Range("J126").Value = 13342.667
$ws.Range("K126").Value = 12741.75
$ws.Range("L126").Value = 40028.001
$ws.Range("M126").Value = -10271.75
$ws.Range("N126").Value = -44968.001

$ws.Range("H132").Value = 25429.355
$ws.Range("I132").Value = 25890.861
$ws.Range("J132").Value = 15507
$ws.Range("K132").Value = 77672.583
$ws.Range("L132").Value = 46521
$ws.Range("M132").Value = -75142.583
$ws.Range("N132").Value = -51581

$ws = $wb.Worksheets.Item(7)
$ws.Range("H20").Value = 25503.75
$ws.Range("I20").Value = 24861.428
$ws.Range("J20").Value = 30000
$ws.Range("K20").Value = 24861.428
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = -24635.428
$ws.Range("N20").Value = -30452

$ws.Range("H35").Value = 2511.1667
$ws.Range("I35").Value = 2500
$ws.Range("J35").Value = 2513.4
$ws.Range("K35").Value = 2500
$ws.Range("L35").Value = 2513.4
$ws.Range("M35").Value = -2164
$ws.Range("N35").Value = -3185.4

$ws.Range("H46").Value = 1965.8422
$ws.Range("I46").Value = 1264.6666
$ws.Range("J46").Value = 2097.3125
$ws.Range("K46").Value = 1264.6666
$ws.Range("L46").Value = 2097.3125
$ws.Range("M46").Value = -1076.6666
$ws.Range("N46").Value = -2473.3125

$ws.Range("H74").Value = 46665
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996

$ws.Range("H77").Value = 46665
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984

$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 50000
$ws.Range("K87").Value = 50000
$ws.Range("M87").Value = -48877

$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 50000
$ws.Range("K90").Value = 150000
$ws.Range("M90").Value = -144384

$ws.Range("H132").Value = 12853
$ws.Range("J132").Value = 19005
$ws.Range("L132").Value = 57015
$ws.Range("N132").Value = -62075

$ws.Range("H136").Value = 11843.909
$ws.Range("I136").Value = 2318.25
$ws.Range("K136").Value = 6954.75
$ws.Range("M136").Value = -4404.75

$ws = $wb.Worksheets.Item(8)
$ws.Range("H41").Value = 157449.5
$ws.Range("J41").Value = 157449.5
$ws.Range("L41").Value = 157449.5
$ws.Range("N41").Value = -158229.5

$ws.Range("H75").Value = 45000
$ws.Range("I75").Value = 45000
$ws.Range("K75").Value = 45000
$ws.Range("M75").Value = -44064

$ws.Range("H78").Value = 45000
$ws.Range("I78").Value = 45000
$ws.Range("K78").Value = 135000
$ws.Range("M78").Value = -130320

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 4089.1904
$ws.Range("I136").Value = 2187.111
$ws.Range("J136").Value = 15501.667
$ws.Range("K136").Value = 6561.333
$ws.Range("L136").Value = 46505.001
$ws.Range("M136").Value = -4011.333
$ws.Range("N136").Value = -51605.001
